$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.257.25'
$ws.Range("E2").Value = '  +0.67%  '

$ws.Range("D3").Value = '2.296.11'
$ws.Range("E3").Value = '  -0.01%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.70'
$ws.Range("E5").Value = '  +1.61%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.55'
$ws.Range("E6").Value = '  -3.32%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.629'
$ws.Range("E7").Value = '  +0.19%  '

$ws.Range("E8").Value = '  +0.21%  '

$ws.Range("E9").Value = '  -0.31%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.56'
$ws.Range("E10").Value = '  -1.99%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0907'
$ws.Range("E11").Value = '  -0.61%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.40'
$ws.Range("E12").Value = '  +1.66%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.106'
$ws.Range("E13").Value = '  +0.58%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.956'
$ws.Range("E14").Value = '  -1.36%  '

$ws.Range("E15").Value = '  -2.04%  '

$ws.Range("D16").Value = '2.643.08'
$ws.Range("E16").Value = '  +0.07%  '

$ws.Range("D17").Value = '2.285.85'
$ws.Range("E17").Value = '  -0.27%  '

$ws.Range("D18").Value = '42.292.66'
$ws.Range("E18").Value = '  +0.55%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.41'
$ws.Range("E19").Value = '  -1.86%  '

$ws.Range("E20").Value = '  +0.34%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.38'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.39'
$ws.Range("E22").Value = '  +0.00%  '

$ws.Range("E23").Value = '  +2.62%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '275.40'
$ws.Range("E24").Value = '  +7.21%  '

$ws.Range("E25").Value = '  -2.45%  '

$ws.Range("E26").Value = '  -0.37%  '

$ws.Range("E27").Value = '  -1.76%  '

$ws.Range("E28").Value = '  +6.02%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.73'
$ws.Range("E29").Value = '  -0.02%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.50'
$ws.Range("E30").Value = '  +5.43%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '165.99'
$ws.Range("E31").Value = '  -0.43%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.07'
$ws.Range("E32").Value = '  +5.30%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0876'
$ws.Range("E33").Value = '  -1.91%  '

$ws.Range("E34").Value = '  +3.94%  '

$ws.Range("E35").Value = '  -8.90%  '

$ws.Range("E36").Value = '  -1.67%  '

$ws.Range("E37").Value = '  -0.26%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0361'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.69'
$ws.Range("E39").Value = '  +2.16%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.75'
$ws.Range("E40").Value = '  -3.26%  '

$ws.Range("E41").Value = '  -0.31%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '96.26'
$ws.Range("E42").Value = '  +0.00%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '69.78'
$ws.Range("E43").Value = '  -2.29%  '

$ws.Range("E44").Value = '  -1.28%  '

$ws.Range("E45").Value = '  -0.21%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '11.95'
$ws.Range("E46").Value = '  -2.97%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '112.60'
$ws.Range("E47").Value = '  -0.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '78.47'
$ws.Range("E48").Value = '  +4.67%  '

$ws.Range("E49").Value = '  -1.57%  '

$ws.Range("E50").Value = '  -0.79%  '

$ws.Range("D51").Value = '1.597.66'
$ws.Range("E51").Value = '  +3.17%  '
